$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Localización"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "ID"
$ws.Range("E1").Value = "Tipo"

# --- Data row (row 2) ----------------------------------------------------
$ws.Range("A2").Value = "Juan Torres Pardo"
$ws.Range("B2").Value = "40.5N30.99W"
$ws.Range("C2").Value = "juan@example.com"
$ws.Range("D2").Value = 123
$ws.Range("E2").Value = 1

# The old sheet had extra columns (F:I) that are no longer part of the
# data set - drop their contents so the used range shrinks back down.
$ws.Range("F1:I2").ClearContents()

# D2 used to hold a date (with a date number format); now it is a plain
# number, so drop the leftover date formatting.
$ws.Range("D2").ClearFormats()

# The "ID" column (old column D, narrower) now carries what used to be the
# "Dirección postal" column's width.
$ws.Columns.Item(4).ColumnWidth = 21.26

# Match the workbook's last recorded selection.
$ws.Range("B2").Select() | Out-Null
